$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 701, pushing existing rows 701-743 down to 704-746
$ws.Rows("701:703").Insert()

# Populate the 3 new rows with new data (columns A-R)

# Row 701
$ws.Cells.Item(701, 1).Value2 = 10
$ws.Cells.Item(701, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(701, 3).Value2 = "La Araucanía"
$ws.Cells.Item(701, 4).Value2 = 45013
$ws.Cells.Item(701, 5).Value2 = 9
$ws.Cells.Item(701, 6).Value2 = 100112032
$ws.Cells.Item(701, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(701, 8).Value2 = "Bola 8"
$ws.Cells.Item(701, 9).Value2 = "Primera"
$ws.Cells.Item(701, 10).Value2 = 40
$ws.Cells.Item(701, 11).Value2 = 14000
$ws.Cells.Item(701, 12).Value2 = 14000
$ws.Cells.Item(701, 13).Value2 = 14000
$ws.Cells.Item(701, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(701, 15).Value2 = "Región del Maule"
$ws.Cells.Item(701, 16).Value2 = 280
$ws.Cells.Item(701, 17).Value2 = 50
$ws.Cells.Item(701, 18).Value2 = "Hortaliza"

# Row 702
$ws.Cells.Item(702, 1).Value2 = 10
$ws.Cells.Item(702, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(702, 3).Value2 = "La Araucanía"
$ws.Cells.Item(702, 4).Value2 = 45013
$ws.Cells.Item(702, 5).Value2 = 9
$ws.Cells.Item(702, 6).Value2 = 100112032
$ws.Cells.Item(702, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(702, 8).Value2 = "Sin especificar"
$ws.Cells.Item(702, 9).Value2 = "Primera"
$ws.Cells.Item(702, 10).Value2 = 50
$ws.Cells.Item(702, 11).Value2 = 12000
$ws.Cells.Item(702, 12).Value2 = 12000
$ws.Cells.Item(702, 13).Value2 = 12000
$ws.Cells.Item(702, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(702, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(702, 16).Value2 = 240
$ws.Cells.Item(702, 17).Value2 = 50
$ws.Cells.Item(702, 18).Value2 = "Hortaliza"

# Row 703
$ws.Cells.Item(703, 1).Value2 = 10
$ws.Cells.Item(703, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(703, 3).Value2 = "La Araucanía"
$ws.Cells.Item(703, 4).Value2 = 45013
$ws.Cells.Item(703, 5).Value2 = 9
$ws.Cells.Item(703, 6).Value2 = 100112032
$ws.Cells.Item(703, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(703, 8).Value2 = "Sin especificar"
$ws.Cells.Item(703, 9).Value2 = "Primera"
$ws.Cells.Item(703, 10).Value2 = 200
$ws.Cells.Item(703, 11).Value2 = 10000
$ws.Cells.Item(703, 12).Value2 = 10000
$ws.Cells.Item(703, 13).Value2 = 10000
$ws.Cells.Item(703, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(703, 15).Value2 = "Región del Maule"
$ws.Cells.Item(703, 16).Value2 = 200
$ws.Cells.Item(703, 17).Value2 = 50
$ws.Cells.Item(703, 18).Value2 = "Hortaliza"
